# Update admin audit IT sheet:
# - Clear the "Status Complete" values in column G (rows 2-7)
# - Move the active selection to G17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the values previously stored in G2:G7 (Status Complete column)
$ws.Range("G2:G7").ClearContents()

# Update the active cell / selection shown in the sheet view
$ws.Range("G17").Select()
